$wb = $excel.ActiveWorkbook

# Work on the "sect1" worksheet (first sheet, tab-selected / active sheet).
$ws = $wb.Worksheets.Item("sect1")
$ws.Activate()

# Rename the header in A1 from "parameter" to "name".
$ws.Range("A1").Value = "name"

# Update the active selection shown in the sheet view.
$ws.Range("C43").Select()

# Reposition the workbook window (matches the bookViews/workbookView diff).
$win = $wb.Windows.Item(1)
$win.Left = 1100
$win.Top = 0
